$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 258; existing rows 258:309 shift down to 259:310.
$ws.Rows("258:258").Insert()

# Populate the newly inserted row 258 with the new data record.
$ws.Cells.Item(258, 1).Value = 3
$ws.Cells.Item(258, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(258, 3).Value = "Coquimbo"
$ws.Cells.Item(258, 4).Value = 44637
$ws.Cells.Item(258, 5).Value = 5
$ws.Cells.Item(258, 6).Value = 100112009
$ws.Cells.Item(258, 7).Value = "Acelga"
$ws.Cells.Item(258, 8).Value = "Sin especificar"
$ws.Cells.Item(258, 9).Value = "Primera"
$ws.Cells.Item(258, 10).Value = 235
$ws.Cells.Item(258, 11).Value = 4000
$ws.Cells.Item(258, 12).Value = 4500
$ws.Cells.Item(258, 13).Value = 4266
$ws.Cells.Item(258, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(258, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(258, 16).Value = 711
$ws.Cells.Item(258, 17).Value = 6
$ws.Cells.Item(258, 18).Value = "Hortaliza"
